$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.791586399078369
$ws.Range("B1").Value = 4.286983966827393
$ws.Range("C1").Value = 1.538731694221497
$ws.Range("D1").Value = 0.8616945147514343
$ws.Range("E1").Value = 0.4655308723449707
